# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Sun Dec 31 01:44:45 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells that are about to be rewritten to Text format first,
# so values such as "317.91" or "1.00" are preserved verbatim (as strings) instead of
# being auto-coerced into numbers (which would drop trailing zeros / introduce float error).
$dCells = @("D2","D3","D5","D6","D9","D10","D11","D12","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Coin name / link swaps (rows reshuffled by the refreshed ranking) ---
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'

# --- Price (column D) updates ---
$ws.Range("D2").Value = '42.421.84'
$ws.Range("D3").Value = '2.293.60'
$ws.Range("D5").Value = '317.91'
$ws.Range("D6").Value = '101.89'
$ws.Range("D9").Value = '0.605'
$ws.Range("D10").Value = '39.23'
$ws.Range("D11").Value = '0.0904'
$ws.Range("D12").Value = '8.31'
$ws.Range("D14").Value = '0.959'
$ws.Range("D15").Value = '15.25'
$ws.Range("D16").Value = '2.643.21'
$ws.Range("D17").Value = '2.305.84'
$ws.Range("D18").Value = '42.384.38'
$ws.Range("D19").Value = '7.40'
$ws.Range("D20").Value = '0.0000105'
$ws.Range("D21").Value = '12.77'
$ws.Range("D22").Value = '73.00'
$ws.Range("D23").Value = '3.57'
$ws.Range("D24").Value = '269.17'
$ws.Range("D25").Value = '2.23'
$ws.Range("D27").Value = '10.83'
$ws.Range("D28").Value = '2.34'
$ws.Range("D29").Value = '22.64'
$ws.Range("D30").Value = '37.23'
$ws.Range("D31").Value = '165.61'
$ws.Range("D32").Value = '6.05'
$ws.Range("D33").Value = '0.0873'
$ws.Range("D34").Value = '0.133'
$ws.Range("D38").Value = '0.0357'
$ws.Range("D39").Value = '3.68'
$ws.Range("D40").Value = '2.75'
$ws.Range("D41").Value = '1.54'
$ws.Range("D42").Value = '69.09'
$ws.Range("D43").Value = '94.48'
$ws.Range("D44").Value = '1.00'
$ws.Range("D45").Value = '0.225'
$ws.Range("D46").Value = '114.99'
$ws.Range("D47").Value = '11.97'
$ws.Range("D48").Value = '79.39'
$ws.Range("D49").Value = '8.92'
$ws.Range("D50").Value = '5.28'
$ws.Range("D51").Value = '1.595.22'

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("E6").Value = '  -4.99%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("E10").Value = '  -2.87%  '
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("E14").Value = '  -2.88%  '
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("E19").Value = '  -3.94%  '
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("E21").Value = '  +28.04%  '
$ws.Range("E22").Value = '  -0.68%  '
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("E24").Value = '  +2.51%  '
$ws.Range("E25").Value = '  -4.69%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("E28").Value = '  +2.76%  '
$ws.Range("E29").Value = '  -1.39%  '
$ws.Range("E30").Value = '  +1.27%  '
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("E32").Value = '  +1.54%  '
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("E34").Value = '  +1.98%  '
$ws.Range("E35").Value = '  -8.92%  '
$ws.Range("E36").Value = '  -4.89%  '
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("E40").Value = '  -6.72%  '
$ws.Range("E41").Value = '  +2.49%  '
$ws.Range("E42").Value = '  -3.54%  '
$ws.Range("E43").Value = '  -8.75%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("E45").Value = '  -1.74%  '
$ws.Range("E46").Value = '  +0.35%  '
$ws.Range("E47").Value = '  -4.76%  '
$ws.Range("E48").Value = '  -2.39%  '
$ws.Range("E49").Value = '  -3.34%  '
$ws.Range("E50").Value = '  -1.74%  '
$ws.Range("E51").Value = '  +2.73%  '

